$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ------------
# Overview sheet: zh-cn/de-de status columns (E/F) for both data rows
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

# zh-cn sheet: Status column (C) for both data rows
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

# de-de sheet: Status column (C) for both data rows
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Narrow the (now shorter) status columns --------------------------------
# Target OOXML column width is ~13.41 characters; this runtime's ColumnWidth
# setter snaps to the nearest 1/6-character pixel grid, so 12.5 is the closest
# settable value (-> stored width 13.3333...).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
